# Apply the two changes described by the commit:
#  1. Bump the cached "datetimeFigureOut" footer field from 8/9/2021 to
#     8/10/2021 on the slide master and every slide layout.
#  2. Relabel the second "Availability Zone" box from "Zone 1" to "Zone 2".

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes, [string]$newDate) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

$newDate = "8/10/2021"

# Slide master footer date field.
Update-DatePlaceholder $p.SlideMaster.Shapes $newDate

# Every slide layout's footer date field.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes $newDate
}

# Relabel the right-hand Availability Zone box (Rectangle 93) from
# "Zone 1" to "Zone 2", leaving the left-hand box untouched.
$slide = $p.Slides.Item(1)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.Name -eq "Rectangle 93") {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq "Availability Zone 1") {
            $lastChar = $tr.Characters($tr.Length, 1)
            $lastChar.Text = "2"
        }
    }
}

Write-Output "edit applied"
